# T2D: Refactoring and improvements - add WAL / GUI menu helper function
# descriptors to the Sheet1 function list and register the new "GUI"
# category on the Category sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cat = $wb.Worksheets.Item("Category")

# --- New "GUI" category -----------------------------------------------
$cat.Cells.Item(6, 1).Value = "GUI"

# --- New function rows in Sheet1 (rows 101-105) ------------------------
# Row 101: f.SelectMenu
$ws.Rows.Item(101).RowHeight = 54
$ws.Cells.Item(101, 1).Value = "Сформировать и отобразить меню, вернуть выбранную опцию."
$ws.Cells.Item(101, 2).Value = "GUI"
$ws.Cells.Item(101, 3).Value = 12800
$ws.Cells.Item(101, 4).Value = "Массив с опциями меню"
$ws.Cells.Item(101, 5).Value = "text-id или текст тайтла меню (самая верхняя жирная строка)"
$ws.Cells.Item(101, 6).Value = "text-id или текст с детальным описанием меню"
$ws.Cells.Item(101, 7).Value = "null"
$ws.Cells.Item(101, 8).Value = "null"
$ws.Cells.Item(101, 9).Value = "Выбранное значение"

# Row 102: f.SelectIntValue (percent/absolute)
$ws.Rows.Item(102).RowHeight = 99
$ws.Cells.Item(102, 1).Value = "Выбрать целочисленное значение в процентах или абсолютное"
$ws.Cells.Item(102, 2).Value = "GUI"
$ws.Cells.Item(102, 3).Value = 12801
$ws.Cells.Item(102, 4).Value = "text-id или текст краткого описания выбираемого значения (например: Мин. партия для покупки)"
$ws.Cells.Item(102, 5).Value = "text-id или текст подробного описания значения (Например: Значение минимальной партии для покупки позволяет предотвратить рейсы…)"
$ws.Cells.Item(102, 6).Value = "null"
$ws.Cells.Item(102, 7).Value = "null"
$ws.Cells.Item(102, 8).Value = "null"
$ws.Cells.Item(102, 9).Value = "Положительное - выбор в процентах; Отрицательное - выбор в абсолютных единицах; null - сбросить значение; `$cGUI.Menu.Close - выбрано закрыть терминал; `$cGUI.Menu.Back - выбрана опция возврата в предыдущее меню"

# Row 103: f.AddMenuIntOption
$ws.Rows.Item(103).RowHeight = 69
$ws.Cells.Item(103, 1).Value = "Добавить в меню опцию редактирования целочисленного параметра, который может быть указан в единицах или процентах (MTB, MTS, etc...)"
$ws.Cells.Item(103, 2).Value = "GUI"
$ws.Cells.Item(103, 3).Value = 12802
$ws.Cells.Item(103, 4).Value = "Массив с опциями меню"
$ws.Cells.Item(103, 5).Value = "Текущее значение параметра"
$ws.Cells.Item(103, 6).Value = "text-id краткого наименования опции"
$ws.Cells.Item(103, 7).Value = "текстовый идентификатор опции, который должен быть возвращен в случае выбора в меню"
$ws.Cells.Item(103, 8).Value = "null"
$ws.Cells.Item(103, 9).Value = "null"

# Row 104: f.SelectPrice
$ws.Rows.Item(104).RowHeight = 69
$ws.Cells.Item(104, 1).Value = "Выбрать цену товара"
$ws.Cells.Item(104, 2).Value = "GUI"
$ws.Cells.Item(104, 3).Value = 12803
$ws.Cells.Item(104, 4).Value = "text-id или текст краткого описания выбираемого значения (например: Цена закупки)"
$ws.Cells.Item(104, 5).Value = "Ware"
$ws.Cells.Item(104, 6).Value = "null"
$ws.Cells.Item(104, 7).Value = "null"
$ws.Cells.Item(104, 8).Value = "null"
$ws.Cells.Item(104, 9).Value = "Выбранная цена в пределах min-max цены товара; null - сбросить значение"

# Row 105: f.WareEditMenu
$ws.Rows.Item(105).RowHeight = 69
$ws.Cells.Item(105, 1).Value = "Меню редактирование параметров товара"
$ws.Cells.Item(105, 2).Value = "GUI"
$ws.Cells.Item(105, 3).Value = 12901
$ws.Cells.Item(105, 4).Value = "WareStruct"
$ws.Cells.Item(105, 5).Value = "null"
$ws.Cells.Item(105, 6).Value = "null"
$ws.Cells.Item(105, 7).Value = "null"
$ws.Cells.Item(105, 8).Value = "null"
$ws.Cells.Item(105, 9).Value = "`$cGUI.Menu.Close - выбрана опция закрытия терминала; `$cGUI.Menu.Remove - выбрана опция удаления товара из списка; null - меню закрыто"

# --- Two extra blank trailer rows (111-112), matching the existing blank
#     row style used for rows 106-110 --------------------------------
$ws.Cells.Item(111, 3).Style = $ws.Cells.Item(110, 3).Style
$ws.Cells.Item(112, 3).Style = $ws.Cells.Item(110, 3).Style

# --- Move the view / selection the author left the sheet on -----------
$ws.Application.ActiveWindow.ScrollRow = 98
$ws.Range("I106").Select()
